# compute frequency lr IM
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 15: "m" label plus the integer series 1..9 (G15:O15)
$ws.Range("F15").Value = "m"
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 2
$ws.Range("I15").Value = 3
$ws.Range("J15").Value = 4
$ws.Range("K15").Value = 5
$ws.Range("L15").Value = 6
$ws.Range("M15").Value = 7
$ws.Range("N15").Value = 8
$ws.Range("O15").Value = 9

# New frequency computations (11.3/3 and 14.9/4) added in B19 / B20
$ws.Range("B19").Formula = "=11.3/3"
$ws.Range("B20").Formula = "=14.9/4"

# Restore the view state the workbook was left in: scrolled so row 3 is at
# the top and D21 is the active selection.
$ws.Range("D21").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
